# Update the "想去人数" (F column) counts that changed between the
# previous data pull and the newly generated output.
#
# The same set of row updates needs to be applied to both the "展览"
# sheet and the "全部类型" sheet, since both sheets mirror the same
# underlying data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    2  = 132
    3  = 47
    5  = 83
    7  = 1225
    8  = 1515
    10 = 377
    19 = 1714
    20 = 65
    23 = 657
    26 = 4121
    30 = 1074
    31 = 133
    33 = 490
    35 = 224
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
